# KP-7577: generating worksheet name for filtered/cascading questions.
# Rename the "Translations question" sheet to the generated cascading-question
# sheet-name convention "@@_question".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations question")
$ws.Name = "@@_question"
